$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new optics parts (rows 183-185), leaving row 182 blank as a separator ---

# Row 183: Aimpoint QRP2 Modular Base Mount
$ws.Range("A183").Value = "aimpoint_qrp2_modular_base_mount"
$ws.Range("B183").Value = "Aimpoint QRP2 Modular Base Mount"
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 0.03
$ws.Range("M183").Value = 300
$ws.Range("P183").Value = 11.8

# Row 184: Aimpoint AR-15 Spacer
$ws.Range("A184").Value = "aimpoint_ar15_spacer"
$ws.Range("B184").Value = "Aimpoint AR-15 Spacer"
$ws.Range("C184").Value = -1
$ws.Range("D184").Value = 0.02
$ws.Range("M184").Value = 200

# Row 185: Aimpoint CompM4
$ws.Range("A185").Value = "aimpoint_compm4_sight"
$ws.Range("B185").Value = "Aimpoint CompM4"
$ws.Range("C185").Value = -1
$ws.Range("D185").Value = 0.14
$ws.Range("M185").Value = 1000
$ws.Range("P185").Value = 9.4

# --- Fill the "strength" (N) and "weight formula" (Q) formulas down through the new rows ---
# (this also backfills N169:N181, which had the style but no formula yet)
$ws.Range("N151:N185").Formula = "=C151-D151*20-E151*0.8-F151*0.6-H151*5+I151*10+J151/300"
$ws.Range("Q169:Q185").Formula = "=P169*0.013+0.02"

# --- Tidy up the two lookup columns a bit narrower to make room on screen ---
$ws.Columns("A").ColumnWidth = 39.166666666666664
$ws.Columns("B").ColumnWidth = 43.5

# --- Update the active selection to reflect where editing left off ---
$ws.Range("G182").Select()
